$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (Sending cluster x Target cluster matrix:
# ECs/FAPs/sCs x ECs/FAPs/sCs), Ligand symbol = Spon2, Receptor symbol = Itga4.
$data = @(
        @("ECs", "Spon2", "Itga4", "ECs", 2, 0.6666666666666666, 0.3355466666666667, 1.00664, 0.02388108083384315, 0.02388108083384315, 3, 1, 26.04517333333333, 78.13552, 0.9210237118384171, 0.921023711838417, 8.739371094755557, 78.65433985279999, 0.0219950417122995, 0.0219950417122995),
        @("ECs", "Spon2", "Itga4", "FAPs", 2, 0.6666666666666666, 0.3355466666666667, 1.00664, 0.02388108083384315, 0.02388108083384315, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.01167750336256582, 0.01167750336256582, 0.1108050032088889, 0.99724502888, 0.0002788714017389096, 0.0002788714017389096),
        @("ECs", "Spon2", "Itga4", "sCs", 2, 0.6666666666666666, 0.3355466666666667, 1.00664, 0.02388108083384315, 0.02388108083384315, 3, 1, 1.903109, 5.709327, 0.06729878479901708, 0.06729878479901708, 0.6385818812533333, 5.74723693128, 0.001607167719804741, 0.001607167719804742),
        @("FAPs", "Spon2", "Itga4", "ECs", 3, 1, 12.878362, 38.635086, 0.9165616424823987, 0.9165616424823989, 3, 1, 26.04517333333333, 78.13552, 0.9210237118384171, 0.921023711838417, 335.4191705394134, 3018.77253485472, 0.8441750060878551, 0.8441750060878552),
        @("FAPs", "Spon2", "Itga4", "FAPs", 3, 1, 12.878362, 38.635086, 0.9165616424823987, 0.9165616424823989, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.01167750336256582, 0.01167750336256582, 4.252722749151333, 38.274504742362, 0.01070315166208707, 0.01070315166208707),
        @("FAPs", "Spon2", "Itga4", "sCs", 3, 1, 12.878362, 38.635086, 0.9165616424823987, 0.9165616424823989, 3, 1, 1.903109, 5.709327, 0.06729878479901708, 0.06729878479901708, 24.508926627458, 220.580339647122, 0.06168348473245658, 0.06168348473245659),
        @("sCs", "Spon2", "Itga4", "ECs", 3, 1, 0.8368233333333333, 2.51047, 0.05955727668375805, 0.05955727668375806, 3, 1, 26.04517333333333, 78.13552, 0.9210237118384171, 0.921023711838417, 21.79520876604444, 196.1568788944, 0.05485366403826245, 0.05485366403826246),
        @("sCs", "Spon2", "Itga4", "FAPs", 3, 1, 0.8368233333333333, 2.51047, 0.05955727668375805, 0.05955727668375806, 2, 0.6666666666666666, 0.3302223333333333, 0.990667, 0.01167750336256582, 0.01167750336256582, 0.2763377537211111, 2.48703978349, 0.0006954802987398477, 0.0006954802987398479),
        @("sCs", "Spon2", "Itga4", "sCs", 3, 1, 0.8368233333333333, 2.51047, 0.05955727668375805, 0.05955727668375806, 3, 1, 1.903109, 5.709327, 0.06729878479901708, 0.06729878479901708, 1.592566017076666, 14.33309415369, 0.00400813234675575, 0.004008132346755751)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
